$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Estimación
# ---------------------------------------------------------------------------
$est = $wb.Worksheets.Item("Estimación")

# Row 1 title: merge B1:C1, change text, shrink row height
$est.Range("C1:P1").UnMerge()
$est.Range("B1").Value = "Estimación VIA"
$est.Range("B1:C1").Merge()
$est.Rows.Item(1).RowHeight = 15

# Extend the small header merge I5:J5 -> I5:K5
$est.Range("I5:J5").UnMerge()
$est.Range("I5:K5").Merge()

# Insert a new row for "CAMBIOS" right above the TOTAL row
$est.Rows.Item(12).Insert()
$est.Range("B12").Value = "CAMBIOS"
$est.Range("B12:D12").Merge()

# TOTAL row (now row 13): extend the sums to include the new row 12
$est.Range("E13:F13").Merge()
$est.Range("E13").Formula = "=SUM(E7:F12)"
$est.Range("G13").Formula = "=SUM(G7:G12)"

# Column widths
$est.Columns.Item("B").ColumnWidth = 14.28515625
$est.Columns.Item("C").ColumnWidth = 8.85546875
$est.Columns.Item("D").ColumnWidth = 2.140625
$est.Columns.Item("K").ColumnWidth = 15.85546875

# ---------------------------------------------------------------------------
# Sheet: Capacidad
# ---------------------------------------------------------------------------
$cap = $wb.Worksheets.Item("Capacidad")

$cap.Range("B4").Value = 6
$cap.Range("D4").Value = 1
$cap.Range("H4").Value = 8
$cap.Range("H5").Value = 3
